$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-05 Tuesday" "2025-08-06 Wednesday"

Replace-Text "40÷4=10, 0" "72÷6=12, 0"
Replace-Text "66÷7=9, 3" "90÷5=18, 0"
Replace-Text "54÷3=18, 0" "29÷9=3, 2"
Replace-Text "11÷4=2, 3" "83÷9=9, 2"
Replace-Text "43÷3=14, 1" "12÷2=6, 0"

Replace-Text "72÷7=10, 2" "35÷6=5, 5"
Replace-Text "75÷7=10, 5" "81÷4=20, 1"
Replace-Text "68÷4=17, 0" "26÷2=13, 0"
Replace-Text "34÷8=4, 2" "55÷8=6, 7"
Replace-Text "18÷8=2, 2" "91÷5=18, 1"

Replace-Text "44÷2=22, 0" "19÷5=3, 4"
Replace-Text "35÷3=11, 2" "60÷8=7, 4"
Replace-Text "67÷2=33, 1" "88÷3=29, 1"
Replace-Text "18÷3=6, 0" "55÷8=6, 7"
Replace-Text "99÷2=49, 1" "11÷7=1, 4"

Replace-Text "88÷9=9, 7" "52÷5=10, 2"
Replace-Text "91÷4=22, 3" "26÷9=2, 8"
Replace-Text "62÷7=8, 6" "68÷6=11, 2"
Replace-Text "32÷7=4, 4" "84÷8=10, 4"
Replace-Text "11÷2=5, 1" "97÷8=12, 1"

Replace-Text "59÷2=29, 1" "32÷3=10, 2"
Replace-Text "25÷7=3, 4" "53÷9=5, 8"
Replace-Text "17÷4=4, 1" "86÷2=43, 0"
Replace-Text "11÷6=1, 5" "24÷2=12, 0"
Replace-Text "24÷6=4, 0" "46÷6=7, 4"
